# Insert a new data row at row 301 (pushes existing rows 301-420 down to
# 302-421) and populate it with the new weekly Apio price observation for
# Femacal de La Calera, Coquimbo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(301).Insert()

$ws.Cells.Item(301, 1).Value = 3
$ws.Cells.Item(301, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(301, 3).Value = 'Coquimbo'
$ws.Cells.Item(301, 4).Value = 44755
$ws.Cells.Item(301, 5).Value = 5
$ws.Cells.Item(301, 6).Value = 100112017
$ws.Cells.Item(301, 7).Value = 'Apio'
$ws.Cells.Item(301, 8).Value = 'Americana (o)'
$ws.Cells.Item(301, 9).Value = 'Primera'
$ws.Cells.Item(301, 10).Value = 250
$ws.Cells.Item(301, 11).Value = 10000
$ws.Cells.Item(301, 12).Value = 11000
$ws.Cells.Item(301, 13).Value = 10480
$ws.Cells.Item(301, 14).Value = '$/docena de matas'
$ws.Cells.Item(301, 15).Value = 'Pan de Azúcar'
$ws.Cells.Item(301, 16).Value = 1747
$ws.Cells.Item(301, 17).Value = 6
$ws.Cells.Item(301, 18).Value = 'Hortaliza'
